$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.893.33"
$ws.Range("E2").Value = "  +5.95%  "
$ws.Range("D3").Value = "2.225.90"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.01"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.73"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.86"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "2.555.73"
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.63"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.71"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.797"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.54"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "2.227.75"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("D19").Value = "41.750.88"
$ws.Range("E19").Value = "  +5.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.80"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.67"
$ws.Range("E23").Value = "  +9.80%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.35"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.91"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  +5.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.65"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000258"
$ws.Range("E40").Value = "  +30.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0239"
$ws.Range("E42").Value = "  +5.71%  "
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.60"
$ws.Range("E44").Value = "  +9.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0980"
$ws.Range("E45").Value = "  +7.35%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.66"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "1.467.41"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.49"
$ws.Range("E49").Value = "  -6.66%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  -1.20%  "
